# Rename the "Device" sheet to "Apparatus" and update the wording used on
# that sheet (and in its header cells) from "Device" to "Apparatus",
# matching the commit: 'Change "Device" to "Apparatus" in excel form,
# simulink, function name'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Rename the worksheet tab itself.
$ws.Name = "Apparatus"

# Update the header row wording (B2/C2 before A1 so the shared-string
# table append order matches the authored workbook).
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Mirror the author's final UI state: the "Apparatus" sheet ends up the
# active/selected tab with A2 selected (previously the "Advance" sheet
# was active).
$ws.Activate()
$ws.Range("A2").Select()
